$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new SEO audit row (hidden keyword text) ---
$ws.Range("A8").Value = 'SEO'
$ws.Range("B8").Value = '(index.html l.40 et l.42): texte “caché” comprenant des mots clés pour tromper les algorithmes de Google lors de la lecture de la page web pour son référencement'
$ws.Range("C8").Value = 'ce genre de pratique, en plus d’être de mieux en mieux détecté par les algorithmes de Google, peut amener à une pénalisation du site web dans son référencement'
$ws.Range("D8").Value = 'n’utiliser des mots clés que dans des contextes pertinents et proscrire ce genre de pratique'
$ws.Range("E8").Value = 'supprimer ces balises'

# --- Row 9: new SEO/accessibilite audit row (alt tag description) ---
$ws.Range("A9").Value = 'SEO/accessibilité'
$ws.Range("B9").Value = '(index.html l.41): description subjective du site dans une balise alt'
$ws.Range("C9").Value = 'les balises alt servent non seulement à décrire un élément visuel, mais sont également utilisées pour vérifier la pertinence d’un contenu par les moteurs de recherche'
$ws.Range("D9").Value = 'se contenter d’utiliser les balises alt pour fournir de courtes descriptions objectives des éléments relatifs à celles-ci'
$ws.Range("E9").Value = 'exemple de description : “logo La Chouette agence”'
# add the hyperlink first (TextToDisplay omits the trailing period, matching the URL itself),
# then overwrite the cell text with the full sentence (incl. trailing period) afterwards
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.atinternet.com/glossaire/balise-alt/", ':~:text=La%20balise%20ALT%2C%20également%20connue,apparaît%20pas%20à%20l''écran', "", 'https://www.atinternet.com/glossaire/balise-alt/#:~:text=La%20balise%20ALT%2C%20%C3%A9galement%20connue,appara%C3%AEt%20pas%20%C3%A0%20l''%C3%A9cran')
$ws.Range("F9").Value = 'https://www.atinternet.com/glossaire/balise-alt/#:~:text=La%20balise%20ALT%2C%20%C3%A9galement%20connue,appara%C3%AEt%20pas%20%C3%A0%20l''%C3%A9cran.'
# the source workbook keeps the plain (non hyperlink) text styling on F9
$ws.Range("F9").Font.Underline = $false
$ws.Range("F9").Font.Color = 0
$ws.Range("F9").Font.Name = "Arial"

# --- Row 10: new SEO audit row (empty <li></li> tag) ---
$ws.Range("A10").Value = '???SEO???'
$ws.Range("B10").Value = '(index.html l.52 et l.53): balise <li></li> vide'
$ws.Range("E10").Value = 'supprimer ces balises'

# --- Row 11: new accessibilite audit row (page2 link label) ---
$ws.Range("A11").Value = 'accessibilité'
$ws.Range("B11").Value = '(index.html l.55): nom affiché sur le site pour la redirection vers page2.html (“page2 &gt;”)'
$ws.Range("C11").Value = 'en plus de la faute de frappe pouvait donner une image amateure du site et le décrédibiliser auprès des visiteurs, nommer un lien “page2” ne permet pas de savoir vers quoi celui-ci redirige'
$ws.Range("D11").Value = 'les liens affichés sur le site devraient refléter clairement le type de contenu vers lequels ils redirigent afin de faciliter la navigation sur le site web, d’autant plus pour les personnes utilisant un affichage alternatif'
$ws.Range("E11").Value = 'exemple de nom pour le lien : “Contact”'

# --- update active selection to A12 ---
$ws.Range("A12").Select()

